# Admin Dashboard referred by ... bug
# Append a new Campus Ambassador row (row 3) with the applicant's details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "sahilharpal1234@gmail.com"
$ws.Range("B3").Value = "Sahil Harpal"
$ws.Range("C3").Value = "CA245368"

# The contact number is all digits, so a plain .Value assignment would be
# auto-detected as a number. Build it as text via a formula in a scratch
# cell, then copy/paste-special the *value* into place so it lands as a
# genuine shared-string text cell (matching the other cells' style) instead
# of a numeric cell or a quote-prefixed text cell.
$ws.Range("ZZ1").Formula = '="7276801998"'
$ws.Range("ZZ1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("E3").Value = "Indian Institute of Technology Jodhpur"
